$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$links = @(
    @{ Row = 2;  Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/isa/bp-14sd/";        Text = "BP-14SD" },
    @{ Row = 3;  Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/isa/bp-20sd/";        Text = "BP-20SD" },
    @{ Row = 4;  Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/picmg_1_0/pci-19s/";  Text = "PCI-19S" },
    @{ Row = 5;  Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/isa/bp-10sd/";        Text = "BP-10SD" },
    @{ Row = 6;  Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/picmg_1_0/pci-17sq/"; Text = "PCI-17SQ" },
    @{ Row = 7;  Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/picmg_1_0/px-8s/";    Text = "PX-8S" },
    @{ Row = 8;  Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/pcisa/ip-4sa/";       Text = "IP-4SA" },
    @{ Row = 9;  Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/picmg_1_3/pxe-19s/";  Text = "PXE-19S" },
    @{ Row = 10; Url = "http://new.nnz-ipc.ru//catalogue/komp_yutery_i_komplektuyuwie/passivnye_kross-platy/picmg_1_0/px-20s3/";  Text = "PX-20S3" }
)

foreach ($link in $links) {
    $cell = $ws.Cells.Item($link.Row, 2)
    $cell.Value = $link.Text
    $ws.Hyperlinks.Add($cell, $link.Url)
    $cell.Style = "Normal"
}
